$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update F5: replace the old claim number with a new one, preserving the
#    existing "Text with quote-prefix" cell style (fill + border + numFmt 49
#    + quotePrefix) that the cell already had. Assigning .Value resets the
#    quote-prefix flag, so we re-apply the format from a sibling cell (F4)
#    afterwards.
# ---------------------------------------------------------------------------
$ws.Cells.Item(5, 6).Value = "0420172008483 "

$ws.Range("F4").Copy()
$ws.Range("F5").PasteSpecial(-4122)  # xlPasteFormats

# ---------------------------------------------------------------------------
# 2) Add row 6 - a new claim-number validation row, same Ambiente/URL/Usuario
#    /Contrasenia as row 5 (pre-producción / tcorvetto / silverarrow).
# ---------------------------------------------------------------------------

# Pre-apply the quote-prefixed text format to F6 so the trailing-space,
# all-digit string is kept as text instead of being parsed as a number.
$ws.Range("F4").Copy()
$ws.Range("F6").PasteSpecial(-4122)

$ws.Cells.Item(6, 2).Value = "preproducciongestion.segurossura.com.ar"
$ws.Cells.Item(6, 3).Value = "https://preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do"
$ws.Cells.Item(6, 4).Value = "tcorvetto"
$ws.Cells.Item(6, 5).Value = "silverarrow"
$ws.Cells.Item(6, 6).Value = "1220170301429 "

# Re-apply formatting for the whole row (B:F) and fix up F6's quote-prefix
# that got cleared by the .Value assignment above.
$ws.Range("B5:F5").Copy()
$ws.Range("B6:F6").PasteSpecial(-4122)
$ws.Range("F4").Copy()
$ws.Range("F6").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) Add row 7 - another new claim-number validation row, same data as rows
#    5/6 except for the claim number itself.
# ---------------------------------------------------------------------------
$ws.Range("F4").Copy()
$ws.Range("F7").PasteSpecial(-4122)

$ws.Cells.Item(7, 2).Value = "preproducciongestion.segurossura.com.ar"
$ws.Cells.Item(7, 3).Value = "https://preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do"
$ws.Cells.Item(7, 4).Value = "tcorvetto"
$ws.Cells.Item(7, 5).Value = "silverarrow"
$ws.Cells.Item(7, 6).Value = "1120170200936 "

$ws.Range("B5:F5").Copy()
$ws.Range("B7:F7").PasteSpecial(-4122)
$ws.Range("F4").Copy()
$ws.Range("F7").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4) Update the saved selection to match the author's final cursor position.
# ---------------------------------------------------------------------------
$ws.Range("A6:A7").Select() | Out-Null

$excel.CutCopyMode = $false | Out-Null
